# Update the "log" worksheet: mark the Monday (row 9) entry as completed
# by re-styling the Files_Worked cell (C9), and filling in the Status (D9)
# and Description (E9) cells, then move the active view/selection down to
# the newly completed row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remember C9's original (not-yet-done) formatting before it gets
# overwritten below - it ends up reused on D9.
$ws.Range("C9").Copy()
$ws.Range("D9").PasteSpecial(-4122)   # xlPasteFormats

# C9 keeps its text but picks up the "done" header/fill style (same style
# already used on the completed C8 cell).
$ws.Range("C8").Copy()
$ws.Range("C9").PasteSpecial(-4122)   # xlPasteFormats

# E9 is newly populated, mirroring the formatting of the row above (E8),
# which already shows a completed day's Description.
$ws.Range("E8").Copy()
$ws.Range("E9").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

$ws.Range("D9").Value = "Done"
$ws.Range("E9").Value = "Implemented the category service , supplier service and Product service"

# Scroll the view down a few rows and move the selection to the freshly
# completed Description cell, matching where the author left off editing.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E9").Select()
